$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values in row 13/14 (columns C, D, E)
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0.25

$ws.Range("C14").Value = 0.25
$ws.Range("E14").Value = 0.75

# Row 18: add label and replace formulas with static values
$ws.Range("B18").Value = "Ergebnisreflektion"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Rows.Item(18).RowHeight = 15

# Row 19: replace formulas with static values
$ws.Range("C19").Value = 0.33
$ws.Range("D19").Value = 0.33
$ws.Range("E19").Value = 0.33

# Update selection to B23
$ws.Range("B23").Select()
